$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = [double]"0.9999883732771242"
    "C" = [double]"0.9990763349190175"
    "D" = [double]"0.9999844079329139"
    "E" = [double]"0.9999999798034312"
    "F" = [double]"0.9999962070785023"
    "G" = [double]"1.085303545036088e-05"
    "H" = [double]"0.000862200808886027"
    "I" = [double]"4.843438172917916e-06"
    "J" = [double]"1.689192772863993e-08"
    "K" = [double]"2.430165050323278e-06"
    "L" = [double]"0.0001805656912053873"
    "M" = [double]"0.00329439454989242"
    "N" = [double]"0.9999069862169936"
    "O" = [double]"0.003434643746218193"
    "P" = [double]"64.86213150418904"
    "Q" = [double]"90.45852382642124"
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col" + "$row").Value = $values[$col]
    }
}
